# Weekly price-report update: insert two new daily records for
# "Vega Modelo de Temuco - Cebollín" at the top of the data block
# (rows 108-109), pushing the existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 108; this shifts every
# row from 108 onward down by two (old 108 -> 110, ..., old 222 -> 224)
# and keeps the date-formatted style on column D.
$ws.Rows.Item(108).Resize(2).Insert()

# --- New row 108 --------------------------------------------------
$ws.Cells.Item(108, 1).Value2  = 10
$ws.Cells.Item(108, 2).Value   = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value   = "La Araucanía"
$ws.Cells.Item(108, 4).Value2  = 44484
$ws.Cells.Item(108, 5).Value2  = 9
$ws.Cells.Item(108, 6).Value2  = 100112037
$ws.Cells.Item(108, 7).Value   = "Cebollín"
$ws.Cells.Item(108, 8).Value   = "Sin especificar"
$ws.Cells.Item(108, 9).Value   = "Primera"
$ws.Cells.Item(108, 10).Value2 = 60
$ws.Cells.Item(108, 11).Value2 = 7000
$ws.Cells.Item(108, 12).Value2 = 8000
$ws.Cells.Item(108, 13).Value2 = 7500
$ws.Cells.Item(108, 14).Value  = "`$/docena de paquetes"
$ws.Cells.Item(108, 15).Value  = "Provincia de Cautín"
$ws.Cells.Item(108, 16).Value2 = 625
$ws.Cells.Item(108, 17).Value2 = 12
$ws.Cells.Item(108, 18).Value  = "Hortaliza"

# --- New row 109 --------------------------------------------------
$ws.Cells.Item(109, 1).Value2  = 10
$ws.Cells.Item(109, 2).Value   = "Vega Modelo de Temuco"
$ws.Cells.Item(109, 3).Value   = "La Araucanía"
$ws.Cells.Item(109, 4).Value2  = 44484
$ws.Cells.Item(109, 5).Value2  = 9
$ws.Cells.Item(109, 6).Value2  = 100112037
$ws.Cells.Item(109, 7).Value   = "Cebollín"
$ws.Cells.Item(109, 8).Value   = "Sin especificar"
$ws.Cells.Item(109, 9).Value   = "Primera"
$ws.Cells.Item(109, 10).Value2 = 20
$ws.Cells.Item(109, 11).Value2 = 5000
$ws.Cells.Item(109, 12).Value2 = 5000
$ws.Cells.Item(109, 13).Value2 = 5000
$ws.Cells.Item(109, 14).Value  = "`$/docena de paquetes"
$ws.Cells.Item(109, 15).Value  = "Región de O'Higgins"
$ws.Cells.Item(109, 16).Value2 = 417
$ws.Cells.Item(109, 17).Value2 = 12
$ws.Cells.Item(109, 18).Value  = "Hortaliza"
